# Insert a new row at row 375 (Hortaliza, Feria Lagunitas de Puerto Montt - Betarraga)
# This shifts the existing rows 375:459 down to 376:460 and leaves a fresh row 375
# that inherits the formatting (date style) from the row that used to occupy 375.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("375:375").Insert()

# Populate the new row 375 with the weekly data point. All fields match the
# (now shifted-down) row 376 except for the date (column D), which gets the
# newest weekly reading.
$ws.Range("A375").Value = 4
$ws.Range("B375").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C375").Value = "Los Lagos"
$ws.Range("D375").Value = 45015
$ws.Range("E375").Value = 10
$ws.Range("F375").Value = 100114014
$ws.Range("G375").Value = "Betarraga"
$ws.Range("H375").Value = "Sin especificar"
$ws.Range("I375").Value = "Primera"
$ws.Range("J375").Value = 500
$ws.Range("K375").Value = 1100
$ws.Range("L375").Value = 1200
$ws.Range("M375").Value = 1150
$ws.Range("N375").Value = '$/paquete 5 unidades'
$ws.Range("O375").Value = "Provincia de Cautín"
$ws.Range("P375").Value = 230
$ws.Range("Q375").Value = 5
$ws.Range("R375").Value = "Hortaliza"
